$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45179 = 2023-09-10).
# This automatic update bumps it by one day (45180 = 2023-09-11) for every
# data row (rows 2 through 82).
$ws.Range("C2:C82").Value = 45180
